$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@("Bitcoin", "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc", "26.279.08", "  +0.64%  ")
    ,@("Ethereum", "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth", "1.663.10", "  +0.62%  ")
    ,@("TetherUSD", "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt", "1.009", "  +0.68%  ")
    ,@("BNB", "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb", "218.48", "  +0.09%  ")
    ,@("XRP", "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp", "0.5315", "  +0.41%  ")
    ,@("USDC", "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc", "1.010", "  +0.68%  ")
    ,@("Cardano", "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada", "0.2639", "  +1.14%  ")
    ,@("Dogecoin", "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge", "0.06363", "  +0.60%  ")
    ,@("Solana", "https://coinranking.com/coin/zNZHO_Sjf+solana-sol", "20.52", "  +0.62%  ")
    ,@("TRON", "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx", "0.07850", "  +1.31%  ")
    ,@("Polkadot", "https://coinranking.com/coin/25W7FG7om+polkadot-dot", "4.551", "  +1.31%  ")
    ,@("WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "1.675.45", "  +0.58%  ")
    ,@("WrappedliquidstakedEther2.0", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth", "1.892.76", "  +0.71%  ")
    ,@("Polygon", "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic", "0.5532", "  +1.25%  ")
    ,@("ShibaInu", "https://coinranking.com/coin/xz24e0BjL+shibainu-shib", "0.0₅8188", "  +0.76%  ")
    ,@("Litecoin", "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc", "65.67", "  +0.65%  ")
    ,@("WrappedBTC", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc", "26.290.85", "  +0.64%  ")
    ,@("Dai", "https://coinranking.com/coin/MoTuySvg7+dai-dai", "1.009", "  +0.66%  ")
    ,@("Uniswap", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni", "4.661", "  +2.77%  ")
    ,@("BitcoinCash", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch", "192.51", "  -0.41%  ")
    ,@("Avalanche", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax", "10.19", "  +1.48%  ")
    ,@("Chainlink", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "6.053", "  +1.01%  ")
    ,@("BinanceUSD", "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd", "1.011", "  +0.67%  ")
    ,@("Monero", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr", "145.10", "  +3.33%  ")
    ,@("Stellar", "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm", "0.1221", "  -1.39%  ")
    ,@("Cosmos", "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom", "7.234", "  -0.46%  ")
    ,@("EthereumClassic", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", "16.19", "  +0.20%  ")
    ,@("Toncoin", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton", "1.484", "  +3.57%  ")
    ,@("Hedera", "https://coinranking.com/coin/jad286TjB+hedera-hbar", "0.05877", "  -1.02%  ")
    ,@("PancakeSwap", "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake", "1.279", "  +0.34%  ")
    ,@("InternetComputer(DFINITY)", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp", "3.589", "  +2.26%  ")
    ,@("Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "3.302", "  +2.16%  ")
    ,@("LidoDAOToken", "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo", "1.614", "  +4.46%  ")
    ,@("ARBITRUM", "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb", "0.9585", "  +1.41%  ")
    ,@("MXToken", "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx", "2.813", "  +1.89%  ")
    ,@("HuobiToken", "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht", "2.429", "  +0.62%  ")
    ,@("ImmutableX", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx", "0.5808", "  +3.26%  ")
    ,@("VeChain", "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet", "0.01612", "  +0.24%  ")
    ,@("FraxShare", "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs", "5.884", "  +0.65%  ")
    ,@("TrustWalletToken", "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt", "0.8571", "  +1.24%  ")
    ,@("PaxDollar", "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp", "1.009", "  +0.66%  ")
    ,@("Maker", "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr", "1.047.44", "  +3.77%  ")
    ,@("Quant", "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt", "104.32", "  +3.28%  ")
    ,@("RocketPoolETH", "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth", "1.805.73", "  +0.53%  ")
    ,@("Aave", "https://coinranking.com/coin/ixgUfzmLR+aave-aave", "57.30", "  +0.78%  ")
    ,@("BabyDogeCoin", "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge", "0.0₈107", "  +4.05%  ")
    ,@("Frax", "https://coinranking.com/coin/KfWtaeV1W+frax-frax", "1.013", "  +0.91%  ")
    ,@("Mantle", "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt", "0.4373", "  +1.95%  ")
    ,@("EnergySwap", "https://coinranking.com/coin/SbWqqTui-+energyswap-ens", "7.990", "  +3.47%  ")
    ,@("Cronos", "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro", "0.05162", "  +0.21%  ")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $data[$i][0]
    $ws.Cells.Item($row, 3).Value = $data[$i][1]

    $dCell = $ws.Cells.Item($row, 4)
    $dOrigStyle = $dCell.Style
    $dCell.NumberFormat = "@"
    $dCell.Value = $data[$i][2]
    $dCell.Style = $dOrigStyle

    $ws.Cells.Item($row, 5).Value = $data[$i][3]
}
